$wb = $excel.ActiveWorkbook

# ALC row 19
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 8335386
$ws.Range("I19").Value = 20834266
$ws.Range("J19").Value = 2798.6667
$ws.Range("K19").Value = 20834266
$ws.Range("L19").Value = 2798.6667
$ws.Range("M19").Value = -20834091
$ws.Range("N19").Value = -3148.6667

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1500
$ws.Range("I33").Value = 988.7222
$ws.Range("J33").Value = 3800.75
$ws.Range("K33").Value = 988.7222
$ws.Range("L33").Value = 3800.75
$ws.Range("M33").Value = -759.7222
$ws.Range("N33").Value = -4258.75

# ALC row 107
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 1528.3846
$ws.Range("I107").Value = 1528.3846
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1528.3846
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 391.6153999999999
$ws.Range("N107").ClearContents()

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3369040.5
$ws.Range("I132").Value = 1600.8889
$ws.Range("J132").Value = 18522518
$ws.Range("K132").Value = 4802.6667
$ws.Range("L132").Value = 55567554
$ws.Range("M132").Value = -2272.6667
$ws.Range("N132").Value = -55572614

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 16129527
$ws.Range("I135").Value = 449.2
$ws.Range("J135").Value = 83334020
$ws.Range("K135").Value = 4042.8
$ws.Range("L135").Value = 750006180
$ws.Range("M135").Value = -1507.8
$ws.Range("N135").Value = -750011250

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2115.2783
$ws.Range("I138").Value = 1461.04
$ws.Range("J138").Value = 2811.2766
$ws.Range("K138").Value = 4383.12
$ws.Range("L138").Value = 8433.8298
$ws.Range("M138").Value = 756.8800000000001
$ws.Range("N138").Value = -18713.8298

# ALC row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 2065
$ws.Range("I141").Value = 1545
$ws.Range("J141").Value = 3105
$ws.Range("K141").Value = 4635
$ws.Range("L141").Value = 9315
$ws.Range("M141").Value = 545
$ws.Range("N141").Value = -19675

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1743.1
$ws.Range("I32").Value = 1607.6947
$ws.Range("J32").Value = 4315.8
$ws.Range("K32").Value = 1607.6947
$ws.Range("L32").Value = 4315.8
$ws.Range("M32").Value = -1320.6947
$ws.Range("N32").Value = -4889.8

# ARM row 35
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1529
$ws.Range("I35").Value = 1529
$ws.Range("K35").Value = 1529
$ws.Range("M35").Value = -1123

# ARM row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1321.32
$ws.Range("I110").Value = 776
$ws.Range("K110").Value = 776
$ws.Range("M110").Value = 1269

# BSM row 64
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 10729938
$ws.Range("J64").Value = 12821354
$ws.Range("L64").Value = 12821354
$ws.Range("N64").Value = -12821804

# BSM row 67
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H67").Value = 10729938
$ws.Range("J67").Value = 12821354
$ws.Range("L67").Value = 12821354
$ws.Range("N67").Value = -12822914

# BSM row 86
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1963.06
$ws.Range("I86").Value = 1962.3062
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1962.3062
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -839.3062
$ws.Range("N86").Value = -4246

# BSM row 89
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 1963.06
$ws.Range("I89").Value = 1962.3062
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 9811.530999999999
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -4195.530999999999
$ws.Range("N89").Value = -21232

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1925
$ws.Range("I105").Value = 1885.7142
$ws.Range("J105").Value = 1980
$ws.Range("K105").Value = 1885.7142
$ws.Range("L105").Value = 1980
$ws.Range("M105").Value = -138.7141999999999
$ws.Range("N105").Value = -5474

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 945.1
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 15757320
$ws.Range("I134").Value = 25000926
$ws.Range("J134").Value = 2552168.2
$ws.Range("K134").Value = 75002778
$ws.Range("L134").Value = 7656504.600000001
$ws.Range("M134").Value = -75000243
$ws.Range("N134").Value = -7661574.600000001

# CRP row 16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 839.5263
$ws.Range("I16").Value = 843.6429000000001
$ws.Range("J16").Value = 828
$ws.Range("K16").Value = 843.6429000000001
$ws.Range("L16").Value = 828
$ws.Range("M16").Value = -556.6429000000001
$ws.Range("N16").Value = -1402

# CRP row 41
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 9750
$ws.Range("J41").Value = 9000
$ws.Range("L41").Value = 9000
$ws.Range("N41").Value = -9856

# CRP row 105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 8140
$ws.Range("I105").Value = 1755.5555
$ws.Range("J105").Value = 17716.666
$ws.Range("K105").Value = 1755.5555
$ws.Range("L105").Value = 17716.666
$ws.Range("M105").Value = -8.555499999999938
$ws.Range("N105").Value = -21210.666

# CRP row 110
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

# CRP row 113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 839.5263
$ws.Range("I113").Value = 843.6429000000001
$ws.Range("J113").Value = 828
$ws.Range("K113").Value = 843.6429000000001
$ws.Range("L113").Value = 828
$ws.Range("M113").Value = 1326.3571
$ws.Range("N113").Value = -5168

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2455.7334
$ws.Range("I132").Value = 1846.909
$ws.Range("J132").Value = 4130
$ws.Range("K132").Value = 5540.727000000001
$ws.Range("L132").Value = 12390
$ws.Range("M132").Value = -3010.727000000001
$ws.Range("N132").Value = -17450

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 46253.863
$ws.Range("J131").Value = 1035.4286
$ws.Range("L131").Value = 3106.2858
$ws.Range("N131").Value = -13186.2858

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 16927776
$ws.Range("I132").Value = 27513984
$ws.Range("J132").Value = 8266334
$ws.Range("K132").Value = 82541952
$ws.Range("L132").Value = 24799002
$ws.Range("M132").Value = -82539422
$ws.Range("N132").Value = -24804062

# LTW row 61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1527.8572
$ws.Range("I61").Value = 1184.2858
$ws.Range("K61").Value = 1184.2858
$ws.Range("M61").Value = -982.2858000000001

# LTW row 113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 1527.8572
$ws.Range("I113").Value = 1184.2858
$ws.Range("K113").Value = 1184.2858
$ws.Range("M113").Value = 985.7141999999999

# LTW row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1784221.2
$ws.Range("I136").Value = 2263911.5
$ws.Range("J136").Value = 2515
$ws.Range("K136").Value = 6791734.5
$ws.Range("L136").Value = 7545
$ws.Range("M136").Value = -6789184.5
$ws.Range("N136").Value = -12645

# WVR row 107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 15436.9
$ws.Range("I107").Value = 28580
$ws.Range("J107").Value = 8359.846
$ws.Range("K107").Value = 85740
$ws.Range("L107").Value = 25079.538
$ws.Range("M107").Value = -83820
$ws.Range("N107").Value = -28919.538

# WVR row 113
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 294.33334
$ws.Range("I113").Value = 210.22223
$ws.Range("J113").Value = 546.6667
$ws.Range("K113").Value = 630.66669
$ws.Range("L113").Value = 1640.0001
$ws.Range("M113").Value = 1539.33331
$ws.Range("N113").Value = -5980.0001
